# Re-applies the recomputed NATMI TPM values (new expression-matrix re-run)
# to the LR-pairs sheet, plus the corresponding target-cluster relabel for the
# MuSCs / Resolving-Mac rows. Values below are taken verbatim from the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 0.06762866666666667
$ws.Range("H2").Value = 0.202886
$ws.Range("I2").Value = 0.0134153952845566
$ws.Range("J2").Value = 0.0134153952845566
$ws.Range("M2").Value = 1.378421333333333
$ws.Range("N2").Value = 4.135264
$ws.Range("O2").Value = 0.01656231489052403
$ws.Range("P2").Value = 0.01794267551419991
$ws.Range("Q2").Value = 0.09322079687822224
$ws.Range("R2").Value = 0.8389871719040001
$ws.Range("S2").Value = 0.0002221900010836776
$ws.Range("T2").Value = 0.0002407080844855267
$ws.Range("G3").Value = 0.06762866666666667
$ws.Range("H3").Value = 0.202886
$ws.Range("I3").Value = 0.0134153952845566
$ws.Range("J3").Value = 0.0134153952845566
$ws.Range("O3").Value = 0.2170932623988173
$ws.Range("P3").Value = 0.2351865659654651
$ws.Range("Q3").Value = 1.221906904408
$ws.Range("R3").Value = 10.997162139672
$ws.Range("S3").Value = 0.002912391928694102
$ws.Range("T3").Value = 0.00315512074804416
$ws.Range("G4").Value = 0.06762866666666667
$ws.Range("H4").Value = 0.202886
$ws.Range("I4").Value = 0.0134153952845566
$ws.Range("J4").Value = 0.0134153952845566
$ws.Range("M4").Value = 17.58286933333333
$ws.Range("N4").Value = 52.748608
$ws.Range("O4").Value = 0.2112656061941426
$ws.Range("P4").Value = 0.22887321273073
$ws.Range("Q4").Value = 1.189106009187555
$ws.Range("R4").Value = 10.701954082688
$ws.Range("S4").Value = 0.002834211617125892
$ws.Range("T4").Value = 0.003070424618829155
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 0.06762866666666667
$ws.Range("H5").Value = 0.202886
$ws.Range("I5").Value = 0.0134153952845566
$ws.Range("J5").Value = 0.0134153952845566
$ws.Range("M5").Value = 19.2082395
$ws.Range("N5").Value = 38.416479
$ws.Range("O5").Value = 0.2307951156866419
$ws.Range("P5").Value = 0.1666869194070983
$ws.Range("Q5").Value = 1.299027626399
$ws.Range("R5").Value = 7.794165758394
$ws.Range("S5").Value = 0.00309620770668127
$ws.Range("T5").Value = 0.002236170912611253
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("G6").Value = 0.06762866666666667
$ws.Range("H6").Value = 0.202886
$ws.Range("I6").Value = 0.0134153952845566
$ws.Range("J6").Value = 0.0134153952845566
$ws.Range("M6").Value = 26.988955
$ws.Range("N6").Value = 80.966865
$ws.Range("O6").Value = 0.3242837008298742
$ws.Range("P6").Value = 0.3513106263825066
$ws.Range("Q6").Value = 1.825227041376667
$ws.Range("R6").Value = 16.42704337239
$ws.Range("S6").Value = 0.004350394030971657
$ws.Range("T6").Value = 0.004712970920586504
$ws.Range("I7").Value = 0.9827953701592058
$ws.Range("J7").Value = 0.9827953701592059
$ws.Range("M7").Value = 1.378421333333333
$ws.Range("N7").Value = 4.135264
$ws.Range("O7").Value = 0.01656231489052403
$ws.Range("P7").Value = 0.01794267551419991
$ws.Range("Q7").Value = 6.829241004917334
$ws.Range("R7").Value = 61.463169044256
$ws.Range("S7").Value = 0.01627736639352589
$ws.Range("T7").Value = 0.01763397842362462
$ws.Range("I8").Value = 0.9827953701592058
$ws.Range("J8").Value = 0.9827953701592059
$ws.Range("O8").Value = 0.2170932623988173
$ws.Range("P8").Value = 0.2351865659654651
$ws.Range("S8").Value = 0.2133582531783152
$ws.Range("T8").Value = 0.2311402681545018
$ws.Range("I9").Value = 0.9827953701592058
$ws.Range("J9").Value = 0.9827953701592059
$ws.Range("M9").Value = 17.58286933333333
$ws.Range("N9").Value = 52.748608
$ws.Range("O9").Value = 0.2112656061941426
$ws.Range("P9").Value = 0.22887321273073
$ws.Range("Q9").Value = 87.11244474498132
$ws.Range("R9").Value = 784.0120027048318
$ws.Range("S9").Value = 0.2076308596414814
$ws.Range("T9").Value = 0.2249355338252245
$ws.Range("D10").Value = "MuSCs"
$ws.Range("I10").Value = 0.9827953701592058
$ws.Range("J10").Value = 0.9827953701592059
$ws.Range("M10").Value = 19.2082395
$ws.Range("N10").Value = 38.416479
$ws.Range("O10").Value = 0.2307951156866419
$ws.Range("P10").Value = 0.1666869194070983
$ws.Range("Q10").Value = 95.1651673211235
$ws.Range("R10").Value = 570.991003926741
$ws.Range("S10").Value = 0.22682437115219
$ws.Range("T10").Value = 0.1638191326593969
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("I11").Value = 0.9827953701592058
$ws.Range("J11").Value = 0.9827953701592059
$ws.Range("M11").Value = 26.988955
$ws.Range("N11").Value = 80.966865
$ws.Range("O11").Value = 0.3242837008298742
$ws.Range("P11").Value = 0.3513106263825066
$ws.Range("Q11").Value = 133.713889729315
$ws.Range("R11").Value = 1203.425007563835
$ws.Range("S11").Value = 0.3187045197936934
$ws.Range("T11").Value = 0.3452664570964581
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.019102
$ws.Range("H12").Value = 0.057306
$ws.Range("I12").Value = 0.003789234556237495
$ws.Range("J12").Value = 0.003789234556237496
$ws.Range("M12").Value = 1.378421333333333
$ws.Range("N12").Value = 4.135264
$ws.Range("O12").Value = 0.01656231489052403
$ws.Range("P12").Value = 0.01794267551419991
$ws.Range("Q12").Value = 0.02633060430933334
$ws.Range("R12").Value = 0.236975438784
$ws.Range("S12").Value = 0.00006275849591446048
$ws.Range("T12").Value = 0.00006798900608976269
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.019102
$ws.Range("H13").Value = 0.057306
$ws.Range("I13").Value = 0.003789234556237495
$ws.Range("J13").Value = 0.003789234556237496
$ws.Range("O13").Value = 0.2170932623988173
$ws.Range("P13").Value = 0.2351865659654651
$ws.Range("Q13").Value = 0.3451327201680001
$ws.Range("R13").Value = 3.106194481512
$ws.Range("S13").Value = 0.0008226172918079324
$ws.Range("T13").Value = 0.0008911770629191698
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.019102
$ws.Range("H14").Value = 0.057306
$ws.Range("I14").Value = 0.003789234556237495
$ws.Range("J14").Value = 0.003789234556237496
$ws.Range("M14").Value = 17.58286933333333
$ws.Range("N14").Value = 52.748608
$ws.Range("O14").Value = 0.2112656061941426
$ws.Range("P14").Value = 0.22887321273073
$ws.Range("Q14").Value = 0.3358679700053333
$ws.Range("R14").Value = 3.022811730048
$ws.Range("S14").Value = 0.0008005349355353073
$ws.Range("T14").Value = 0.0008672542866763777
$ws.Range("D15").Value = "MuSCs"
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.019102
$ws.Range("H15").Value = 0.057306
$ws.Range("I15").Value = 0.003789234556237495
$ws.Range("J15").Value = 0.003789234556237496
$ws.Range("M15").Value = 19.2082395
$ws.Range("N15").Value = 38.416479
$ws.Range("O15").Value = 0.2307951156866419
$ws.Range("P15").Value = 0.1666869194070983
$ws.Range("Q15").Value = 0.3669157909290001
$ws.Range("R15").Value = 2.201494745574
$ws.Range("S15").Value = 0.0008745368277706539
$ws.Range("T15").Value = 0.0006316158350901515
$ws.Range("D16").Value = "Resolving-Mac"
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.019102
$ws.Range("H16").Value = 0.057306
$ws.Range("I16").Value = 0.003789234556237495
$ws.Range("J16").Value = 0.003789234556237496
$ws.Range("M16").Value = 26.988955
$ws.Range("N16").Value = 80.966865
$ws.Range("O16").Value = 0.3242837008298742
$ws.Range("P16").Value = 0.3513106263825066
$ws.Range("Q16").Value = 0.5155430184100001
$ws.Range("R16").Value = 4.63988716569
$ws.Range("S16").Value = 0.001228787005209141
$ws.Range("T16").Value = 0.001331198365462034
